$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for new columns I (I0) and J (IF)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match header style of existing headers (bold font, thin border, centered/top aligned)
$headerRange = $ws.Range("I1:J1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# Data values for columns I and J, rows 2-80
$data = @(@(6,6), @(6,6), @(9,9), @(7,7), @(6,7), @(8,8), @(7,7), @(7,7), @(7,7), @(7,7), @(7,7), @(8,8), @(9,9), @(8,8), @(8,8), @(6,6), @(7,8), @(7,7), @(7,8), @(7,7), @(8,8), @(8,8), @(8,8), @(9,9), @(7,7), @(6,6), @(8,8), @(7,7), @(7,7), @(8,8), @(1,1), @(8,8), @(6,6), @(10,10), @(7,7), @(7,7), @(7,7), @(9,9), @(7,7), @(8,8), @(7,7), @(8,8), @(7,7), @(8,8), @(8,8), @(7,7), @(7,7), @(7,7), @(7,7), @(8,8), @(8,8), @(7,7), @(7,7), @(9,9), @(7,7), @(10,10), @(8,8), @(8,8), @(8,8), @(8,8), @(7,7), @(7,7), @(8,8), @(7,7), @(7,7), @(7,7), @(8,9), @(9,9), @(9,9), @(9,9), @(9,9), @(8,8), @(9,9), @(7,7), @(5,5), @(7,7), @(4,4), @(5,5), @(3,3))

for ($idx = 0; $idx -lt $data.Count; $idx++) {
    $row = $idx + 2
    $pair = $data[$idx]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
